$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) store numeric/percentage-looking
# values as literal text (inline strings) in the source data, e.g.
# "0.1730" / "-6.27%". Force Text number format on each target cell first
# so Excel does not coerce the assignment into a real number (which would
# drop significant trailing zeros and/or reformat the percentage string).
$textCells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6",
    "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11",
    "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16",
    "D17", "E17", "D18", "E18", "D19", "E19", "E20", "D21", "E21", "D22",
    "E22", "D23", "E23", "D24", "E24", "E25", "D26", "D38", "E38", "D39",
    "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44",
    "E44", "D45", "E45", "D46", "E46", "E47", "D48", "E48", "D49", "E49",
    "D50", "E50"
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '293.62'
$ws.Range("E2").Value = '-6.27%'
$ws.Range("D3").Value = '40.74'
$ws.Range("E3").Value = '-0.33%'
$ws.Range("D4").Value = '5.025'
$ws.Range("E4").Value = '-2.18%'
$ws.Range("D5").Value = '0.07396'
$ws.Range("E5").Value = '-3.06%'
$ws.Range("D6").Value = '4.279'
$ws.Range("E6").Value = '-1.25%'
$ws.Range("D7").Value = '1.554'
$ws.Range("E7").Value = '-8.01%'
$ws.Range("D8").Value = '0.9243'
$ws.Range("E8").Value = '-1.04%'
$ws.Range("D9").Value = '2.349'
$ws.Range("E9").Value = '-3.13%'
$ws.Range("D10").Value = '0.1149'
$ws.Range("E10").Value = '-8.41%'
$ws.Range("D11").Value = '0.1730'
$ws.Range("E11").Value = '-5.25%'
$ws.Range("D12").Value = '0.08718'
$ws.Range("E12").Value = '-3.57%'
$ws.Range("D13").Value = '0.04177'
$ws.Range("E13").Value = '0.26%'
$ws.Range("D14").Value = '0.1053'
$ws.Range("E14").Value = '-0.36%'
$ws.Range("D15").Value = '0.001261'
$ws.Range("E15").Value = '-0.46%'
$ws.Range("D16").Value = '0.005961'
$ws.Range("E16").Value = '2.34%'
$ws.Range("D17").Value = '3.416'
$ws.Range("E17").Value = '1.51%'
$ws.Range("D18").Value = '0.3284'
$ws.Range("E18").Value = '-2.25%'
$ws.Range("D19").Value = '7.689'
$ws.Range("E19").Value = '-8.87%'
$ws.Range("E20").Value = '2.28%'
$ws.Range("D21").Value = '0.2878'
$ws.Range("E21").Value = '4.95%'
$ws.Range("D22").Value = '0.03865'
$ws.Range("E22").Value = '-4.35%'
$ws.Range("D23").Value = '0.001258'
$ws.Range("E23").Value = '-0.58%'
$ws.Range("D24").Value = '0.003870'
$ws.Range("E24").Value = '-4.48%'
$ws.Range("E25").Value = '0.27%'
$ws.Range("D26").Value = '0.0003715'
$ws.Range("D38").Value = '0.02339'
$ws.Range("E38").Value = '-5.79%'
$ws.Range("D39").Value = '0.05024'
$ws.Range("E39").Value = '-3.30%'
$ws.Range("D40").Value = '0.005469'
$ws.Range("E40").Value = '152.31%'
$ws.Range("D41").Value = '0.007685'
$ws.Range("E41").Value = '-1.38%'
$ws.Range("D42").Value = '0.1286'
$ws.Range("E42").Value = '-1.04%'
$ws.Range("D43").Value = '0.007326'
$ws.Range("E43").Value = '-0.55%'
$ws.Range("D44").Value = '0.007104'
$ws.Range("E44").Value = '-13.22%'
$ws.Range("D45").Value = '0.3161'
$ws.Range("E45").Value = '1.13%'
$ws.Range("D46").Value = '0.00006366'
$ws.Range("E46").Value = '-4.37%'
$ws.Range("E47").Value = '-0.55%'
$ws.Range("D48").Value = '0.01692'
$ws.Range("E48").Value = '-93.43%'
$ws.Range("D49").Value = '0.00002096'
$ws.Range("E49").Value = '-0.55%'
$ws.Range("D50").Value = '0.0001996'
$ws.Range("E50").Value = '-0.55%'

# Columns B (Coin) and C (Link) are plain text; no numeric coercion risk,
# so they can be assigned directly.
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
